# "add status in rm and pt" - update raw material sheet data:
# - sap_code value "ASTB20_20" -> "ASTB20_32" (shared text, used across rows)
# - C2 ("container") value changes to "container"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rm")

# Replace the sap_code text everywhere it appears (shared string update)
$ws.Cells.Replace("ASTB20_20", "ASTB20_32")

# Update the container cell in row 2
$ws.Range("C2").Value = "container"

# Reflect the new selection/active cell
$ws.Range("C2").Select()
